$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ELSAN"
$ws.Range("B2").Value = "58 bis Rue La Boétie, 75008 Paris, France"
$ws.Range("C2").Value = "elsan.care"
$ws.Range("D2").Value = "+33 1 58 56 16 80"
$ws.Range("E2").Value = 28
$ws.Range("F2").Value = 3.5

# Row 3
$ws.Range("A3").Value = "Dentego"
$ws.Range("B3").Value = "19 Rue de Passy, 75016 Paris, France"
$ws.Range("C3").Value = "dentego.fr"
$ws.Range("D3").Value = "+33 1 88 88 09 09"
$ws.Range("E3").Value = 104
$ws.Range("F3").Value = 4

# Row 4
$ws.Range("A4").Value = "Paris Dental Studios - MARAIS"
$ws.Range("B4").Value = "28 Rue Meslay, 75003 Paris, France"
$ws.Range("C4").Value = "parisdentalstudios.com"
$ws.Range("D4").Value = "+33 9 52 34 01 45"
$ws.Range("E4").Value = 204
$ws.Range("F4").Value = 4.8

# Row 5
$ws.Range("A5").Value = "Clinadent - Centre dentaire Paris 16 Victor Hugo"
$ws.Range("B5").Value = "3 Pl. Victor Hugo, 75016 Paris, France"
$ws.Range("C5").Value = "centre-dentaire-paris16.fr"
$ws.Range("D5").Value = "+33 1 42 25 40 79"
$ws.Range("E5").Value = 601
$ws.Range("F5").Value = 4.6

# Row 6
$ws.Range("A6").Value = "Dentego"
$ws.Range("B6").Value = "111 Av. du Général Leclerc, 75014 Paris, France"
$ws.Range("C6").Value = "dentego.fr"
$ws.Range("D6").Value = "+33 1 40 43 41 00"
$ws.Range("E6").Value = 307
$ws.Range("F6").Value = 3.6

# Row 7
$ws.Range("A7").Value = "Place dentaire - Centre dentaire Paris Nation Saint Antoine"
$ws.Range("B7").Value = "238 bis Rue du Faubourg Saint-Antoine, 75012 Paris, France"
$ws.Range("C7").Value = "centre-dentaire-nation-saint-antoine.fr"
$ws.Range("D7").Value = "+33 1 42 55 55 42"
$ws.Range("E7").Value = 82
$ws.Range("F7").Value = 3.9
